$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-09 10:41:20"
$wsZh.Range("H4").Value = "2016-03-09 10:42:15"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-09 10:41:31"
$wsDe.Range("H4").Value = "2016-03-09 10:42:32"
